$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, pushing existing rows 24..111 down to 25..112
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new record's data
$ws.Range("A24").Value = 10
$ws.Range("B24").Value = 'Vega Modelo de Temuco'
$ws.Range("C24").Value = 'La Araucanía'
$ws.Range("D24").Value = 45133
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = 'Fruta'
$ws.Range("G24").Value = 100108
$ws.Range("H24").Value = 'Tropicales y subtropicales'
$ws.Range("I24").Value = 100108004
$ws.Range("J24").Value = 'Papaya'
$ws.Range("K24").Value = 'Cultivar IV Región'
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 25000
$ws.Range("O24").Value = 25000
$ws.Range("P24").Value = 25000
$ws.Range("Q24").Value = '$/bandeja 10 kilos'
$ws.Range("R24").Value = 'Provincia del Elquí'
$ws.Range("S24").Value = 2500
$ws.Range("T24").Value = 10
